$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 45877.37517244972
$ws.Range("D10").Value = 14.83
$ws.Range("E10").Value = 92.65000000000001
$ws.Range("F10").Value = 163.48
$ws.Range("G10").Value = 5.36
$ws.Range("H10").Value = "SE"
$ws.Range("J10").Value = "09:00:14"
